$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in today's labeling progress numbers for row 6 (E6:I6)
$ws.Range("E6").Value = 370
$ws.Range("F6").Value = 355
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0

# Update the active selection to B4
$ws.Range("B4").Select()
